$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.714.98'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.623.44'
$ws.Range("E3").Value = '  -1.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.991'
$ws.Range("E4").Value = '  -0.90%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.53'
$ws.Range("E5").Value = '  -1.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("E7").Value = '  -0.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.17'
$ws.Range("E8").Value = '  -2.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.256'
$ws.Range("E9").Value = '  -3.42%  '

$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0877'
$ws.Range("E11").Value = '  -0.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.852.98'
$ws.Range("E12").Value = '  -1.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.624.38'
$ws.Range("E13").Value = '  -1.34%  '

$ws.Range("E14").Value = '  -1.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.12'
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.718.25'
$ws.Range("E17").Value = '  -0.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.17'
$ws.Range("E18").Value = '  -0.63%  '

$ws.Range("E19").Value = '  -0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("E20").Value = '  -2.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.991'
$ws.Range("E21").Value = '  -0.92%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.34'
$ws.Range("E22").Value = '  -1.78%  '

$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.28'
$ws.Range("E23").Value = '  -4.76%  '

$ws.Range("E24").Value = '  -5.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.20'
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").Value = '  -0.73%  '

$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.58'
$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.992'
$ws.Range("E29").Value = '  -0.87%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.17'
$ws.Range("E30").Value = '  -1.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0479'
$ws.Range("E31").Value = '  -1.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("E32").Value = '  +1.22%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.07'
$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.390.74'
$ws.Range("E34").Value = '  -3.90%  '

$ws.Range("E35").Value = '  -1.04%  '

$ws.Range("E36").Value = '  +8.62%  '

$ws.Range("E37").Value = '  +0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0170'
$ws.Range("E38").Value = '  +0.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.555'
$ws.Range("E39").Value = '  -1.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.866'
$ws.Range("E40").Value = '  -3.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.01'
$ws.Range("E41").Value = '  -1.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.991'
$ws.Range("E42").Value = '  -0.94%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.84'
$ws.Range("E43").Value = '  +0.86%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.04'
$ws.Range("E44").Value = '  -4.60%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.47'
$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.18'
$ws.Range("E46").Value = '  -1.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.764.21'
$ws.Range("E47").Value = '  -1.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.79'
$ws.Range("E48").Value = '  -1.57%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  -2.63%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0996'
$ws.Range("E50").Value = '  -1.30%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0506'
$ws.Range("E51").Value = '  -0.33%  '
